# Commit working 3D Mapping module
#
# Applies the edits to the "G10 element characteristics" sheet:
#  - change three input cells (C26, C27, C28)
#  - change the formulas in C32, C33 and C43 (order/terms tweak)
#  - fill in the previously-empty ratio rows 49-52 (mirrors rows 36-39)
#  - fix the unit label in D39 (mm2 -> m2)
#  - move the selection / activate the "G10 element characteristics" tab
#    (this also flips tabSelected off "Heat Deposition" and onto this sheet,
#    and sets the workbook's activeTab)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("G10 element characteristics")

# Make this the active sheet -- this is what flips <sheetView tabSelected="1">
# from "Heat Deposition" onto "G10 element characteristics", and sets
# workbookView/@activeTab="3" in xl/workbook.xml.
$ws.Activate()

# --- input cells ---
$ws.Range("C26").Value = 10
$ws.Range("C27").Value = 32
$ws.Range("C28").Value = 8

# --- formula tweaks ---
$ws.Range("C32").Formula = "=(2*C24+2*C25)"
$ws.Range("C33").Formula = "=(C32+1)*(C28-1)*C29"
$ws.Range("C43").Formula = "=(C27)*(2*C21+2*C22)"

# --- unit fix ---
$ws.Range("D39").Value = "m2"

# --- new rows 49-52 (mirrors the C36:D39 block above) ---
$ws.Range("C49").Formula = "=C44/C48"
$ws.Range("D49").Value = "mm3"

$ws.Range("C50").Formula = "=C37"
$ws.Range("D50").Value = "mm"

$ws.Range("C51").Formula = "=C49/C50"
$ws.Range("D51").Value = "mm2"

$ws.Range("C52").Formula = "=C51*0.000001"
$ws.Range("C52").Font.Bold = $true
$ws.Range("D52").Value = "m2"
$ws.Range("D52").Font.Bold = $true

# --- final selection on this sheet ---
$ws.Range("O32").Select()
